$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.642.91'
$ws.Range('E2').Value = '  +1.67%  '
$ws.Range('D3').Value = '2.844.57'
$ws.Range('E3').Value = '  +3.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '359.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '116.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.550'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.03%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.604'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.98'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0867'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.92'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.68%  '
$ws.Range('E15').Value = '  +2.93%  '
$ws.Range('D16').Value = '2.833.85'
$ws.Range('E16').Value = '  +2.66%  '
$ws.Range('E17').Value = '  +2.31%  '
$ws.Range('D18').Value = '52.550.36'
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('B19').Value = 'ImmutableX'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.08%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.32'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.72'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = '0.0₃0989'
$ws.Range('E22').Value = '  +2.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '272.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.85'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.47%  '
$ws.Range('E26').Value = '  +1.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.37'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.76%  '
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.71'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '51.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0456'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +32.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.87'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0842'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.13'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.95'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.29'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.96'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.65'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.90'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.24%  '
$ws.Range('E43').Value = '  +2.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.29'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.07%  '
$ws.Range('D47').Value = '2.075.89'
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.29'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.975'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +11.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.43%  '
